$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 9012
$ws.Range("F3").Value = 1977
$ws.Range("F4").Value = 6616
$ws.Range("F5").Value = 177
$ws.Range("F6").Value = 2142
$ws.Range("F7").Value = 601
$ws.Range("F8").Value = 80
$ws.Range("F13").Value = 6
$ws.Range("F14").Value = 83
$ws.Range("F16").Value = 8913
$ws.Range("F17").Value = 169
$ws.Range("F21").Value = 1843
$ws.Range("F24").Value = 5
$ws.Range("F25").Value = 90
$ws.Range("F27").Value = 201
$ws.Range("F28").Value = 1040
$ws.Range("F29").Value = 15
$ws.Range("F30").Value = 72
$ws.Range("F31").Value = 554
$ws.Range("F32").Value = 29
$ws.Range("F33").Value = 26
$ws.Range("F34").Value = 545
$ws.Range("F35").Value = 2329
$ws.Range("F36").Value = 877
$ws.Range("F37").Value = 542
$ws.Range("F41").Value = 295
$ws.Range("F42").Value = 182
$ws.Range("F44").Value = 1062
$ws.Range("F46").Value = 18
$ws.Range("F47").Value = 81
$ws.Range("F48").Value = 7
$ws.Range("F49").Value = 3999

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 408
$ws.Range("F9").Value = 4
$ws.Range("F14").Value = 16
$ws.Range("F16").Value = 29
$ws.Range("F18").Value = 98

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 2342
$ws.Range("F3").Value = 724
$ws.Range("F4").Value = 334
$ws.Range("F5").Value = 10

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 2342
$ws.Range("F3").Value = 9012
$ws.Range("F4").Value = 334
$ws.Range("F5").Value = 1977
$ws.Range("F6").Value = 6616
$ws.Range("F7").Value = 177
$ws.Range("F8").Value = 2142
$ws.Range("F10").Value = 4
$ws.Range("F11").Value = 601
$ws.Range("F18").Value = 83
$ws.Range("F19").Value = 8914
$ws.Range("F20").Value = 169
$ws.Range("F23").Value = 1843
$ws.Range("F26").Value = 90
$ws.Range("F27").Value = 201
$ws.Range("F28").Value = 1040
$ws.Range("F29").Value = 15
$ws.Range("F30").Value = 72
$ws.Range("F32").Value = 554
$ws.Range("F33").Value = 29
$ws.Range("F34").Value = 27
$ws.Range("F35").Value = 545
$ws.Range("F36").Value = 2329
$ws.Range("F37").Value = 877
$ws.Range("F38").Value = 16
$ws.Range("F40").Value = 542
$ws.Range("F41").Value = 295
$ws.Range("F42").Value = 182
$ws.Range("F43").Value = 81
$ws.Range("F44").Value = 3999
$ws.Range("F45").Value = 68
